# Trade #99 closed at 2026-02-17 09:17:45 - unknown UNKNOWN +0.000%
#
# Updates summary/status rollups for the new closed trade and appends the
# trade's row (Trade # 99, zero-based row 100) to both the "All Trades" and
# "MarketMaking" detail sheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Summary sheet: bump Total Trades and recompute Win Rate %
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B6").Value = 99      # Total Trades
$summary.Range("B9").Value = 41.41   # Win Rate %

# ---------------------------------------------------------------------
# 2) Strategy Status sheet: same rollups for the MarketMaking strategy row
# ---------------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("D4").Value = 99       # Trades
$status.Range("G4").Value = 41.41    # Win Rate %

# ---------------------------------------------------------------------
# 3) Append the new trade row (row 100) to a detail sheet
# ---------------------------------------------------------------------
function Add-TradeRow {
    param($ws)

    $row = 100

    $ws.Cells.Item($row, 1).Value = 99                 # A: Trade #

    # Force date/time-looking strings to stay text, matching the rest of
    # the column (Excel would otherwise coerce "2026-02-17" into a date
    # serial number).
    $ws.Cells.Item($row, 2).NumberFormat = "@"
    $ws.Cells.Item($row, 2).Value = "2026-02-17"       # B: Date
    $ws.Cells.Item($row, 3).NumberFormat = "@"
    $ws.Cells.Item($row, 3).Value = "09:17:39"         # C: Time

    $ws.Cells.Item($row, 4).Value = "MarketMaking"     # D: Strategy
    $ws.Cells.Item($row, 5).Value = "UP"               # E: Side
    $ws.Cells.Item($row, 6).Value = 0.14                # F: Entry Price
    $ws.Cells.Item($row, 7).Value = 0.14                # G: Exit Price
    $ws.Cells.Item($row, 8).Value = "CLOSED"           # H: Status
    $ws.Cells.Item($row, 9).Value = 0                   # I: P&L %
    $ws.Cells.Item($row, 10).Value = 0                  # J: P&L $
    $ws.Cells.Item($row, 11).Value = 100.11             # K: Capital After
    $ws.Cells.Item($row, 12).Value = 0                  # L: Entry Slippage (bps)
    $ws.Cells.Item($row, 13).Value = 0                  # M: Exit Slippage (bps)
    $ws.Cells.Item($row, 14).Value = 0.6                # N: Confidence
    $ws.Cells.Item($row, 15).Value = "Normal spread capture: 19600 bps"  # O: Entry Reason
    $ws.Cells.Item($row, 16).Value = "early_exit"       # P: Exit Reason
    $ws.Cells.Item($row, 17).Value = 0.15               # Q: Duration (min)
}

Add-TradeRow $wb.Worksheets.Item("All Trades")
Add-TradeRow $wb.Worksheets.Item("MarketMaking")
